$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.765.39"
$ws.Range("E2").Value = "  -1.46%  "

$ws.Range("D3").Value = "2.445.88"
$ws.Range("E3").Value = "  +0.26%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.25"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.34%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.21"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.97%  "

$ws.Range("E7").Value = "  +0.15%  "

$ws.Range("E8").Value = "  -0.43%  "

$ws.Range("D9").Value = "2.437.85"
$ws.Range("E9").Value = "  +0.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.109"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.43%  "

$ws.Range("E11").Value = "  +2.50%  "

$ws.Range("E12").Value = "  -1.12%  "

$ws.Range("E13").Value = "  -2.63%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.79"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.60%  "

$ws.Range("D15").Value = "2.889.74"

$ws.Range("E16").Value = "  -1.73%  "

$ws.Range("D17").Value = "61.646.46"
$ws.Range("E17").Value = "  -1.28%  "

$ws.Range("D18").Value = "2.412.25"
$ws.Range("E18").Value = "  -1.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.60"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.16"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.55%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.21"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.98%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.53"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +9.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.07"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.77%  "

$ws.Range("E24").Value = "  -0.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.93"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.31%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "64.89"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.50%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.12"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.89%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "579.62"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -8.94%  "

$ws.Range("D29").Value = "2.567.12"
$ws.Range("E29").Value = "  +0.08%  "

$ws.Range("E30").Value = "  -0.11%  "

$ws.Range("D31").Value = "0.0₃0925"
$ws.Range("E31").Value = "  -3.92%  "

$ws.Range("E32").Value = "  -2.68%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.37"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -5.37%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.87"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.58%  "

$ws.Range("E35").Value = "  -3.78%  "

$ws.Range("E36").Value = "  +0.12%  "

$ws.Range("E37").Value = "  -5.03%  "

$ws.Range("E38").Value = "  -1.50%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "151.61"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.75%  "

$ws.Range("E40").Value = "  -4.59%  "

$ws.Range("E41").Value = "  -1.55%  "

$ws.Range("E42").Value = "  -3.39%  "

$ws.Range("E43").Value = "  +0.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.03"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.07%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.67"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.44%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.35"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -5.59%  "

$ws.Range("D47").Value = "0.0₆0277"
$ws.Range("E47").Value = "  +17.81%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "140.78"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.47%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.56"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.93%  "

$ws.Range("E50").Value = "  -0.19%  "

$ws.Range("E51").Value = "  -2.94%  "
